# Update the division problems in the single table of the worksheet.
# The document contains one 5-column table; every 4th row (1, 5, 9, 13, 17)
# holds the actual "a÷b=" expressions, the rows in between are blank
# spacer rows. We update each populated cell in place using
# Cell.Range.Text so the existing run formatting (font/size) is kept.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1,1).Range.Text = "79÷5="
$t.Cell(1,2).Range.Text = "94÷3="
$t.Cell(1,3).Range.Text = "23÷2="
$t.Cell(1,4).Range.Text = "11÷7="
$t.Cell(1,5).Range.Text = "64÷7="

# Row 5
$t.Cell(5,1).Range.Text = "67÷9="
$t.Cell(5,2).Range.Text = "91÷5="
$t.Cell(5,3).Range.Text = "38÷8="
$t.Cell(5,4).Range.Text = "95÷7="
$t.Cell(5,5).Range.Text = "52÷2="

# Row 9
$t.Cell(9,1).Range.Text = "38÷5="
$t.Cell(9,2).Range.Text = "22÷4="
$t.Cell(9,3).Range.Text = "44÷4="
$t.Cell(9,4).Range.Text = "73÷7="
$t.Cell(9,5).Range.Text = "85÷9="

# Row 13
$t.Cell(13,1).Range.Text = "22÷3="
$t.Cell(13,2).Range.Text = "83÷5="
$t.Cell(13,3).Range.Text = "66÷3="
$t.Cell(13,4).Range.Text = "67÷8="
$t.Cell(13,5).Range.Text = "58÷3="

# Row 17
$t.Cell(17,1).Range.Text = "75÷5="
$t.Cell(17,2).Range.Text = "10÷9="
$t.Cell(17,3).Range.Text = "51÷9="
$t.Cell(17,4).Range.Text = "67÷5="
$t.Cell(17,5).Range.Text = "79÷2="
